$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Due dates shifted (rows 23-29). Values below are the final, already-sorted
# (ascending by due date) contents for columns B (assessment name), C (due
# date serial) and D (due time).
# ---------------------------------------------------------------------------

$ws.Range("B23").Value2 = "Quiz 06"
$ws.Range("C23").Value2 = 45343
$ws.Range("D23").Value2 = "in class"

$ws.Range("B24").Value2 = "Module 7 Programming Problems"
$ws.Range("C24").Value2 = 45349
$ws.Range("D24").Value2 = "7pm"

$ws.Range("B25").Value2 = "Short Programming Project 5"
$ws.Range("C25").Value2 = 45349
$ws.Range("D25").Value2 = "7pm"

$ws.Range("B26").Value2 = "Quiz 07"
$ws.Range("C26").Value2 = 45350
$ws.Range("D26").Value2 = "in class"

$ws.Range("B27").Value2 = "Programming Project 5"
$ws.Range("C27").Value2 = 45351
$ws.Range("D27").Value2 = "7pm"

$ws.Range("B28").Value2 = "Short Programming Project 6"
$ws.Range("C28").Value2 = 45356
$ws.Range("D28").Value2 = "7pm"

$ws.Range("B29").Value2 = "Programming Project 6"
$ws.Range("C29").Value2 = 45358
$ws.Range("D29").Value2 = "7pm"

# ---------------------------------------------------------------------------
# The "next assessment due" highlight (bold black font on the due-date cell)
# moves from C29 to C27, since Programming Project 5 (row 27) is now the
# earliest still-upcoming due date.
# ---------------------------------------------------------------------------

$ws.Range("C29").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("C28").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null   # xlPasteFormats, restores plain date style

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# A stray formatted-but-empty cell in column F also moved: it used to sit
# next to the highlighted row (F29) and a trailing row (F32); now it sits
# next to the new highlighted row (F27) and the row after it (F30).
# ---------------------------------------------------------------------------

$ws.Range("F29").Copy() | Out-Null
$ws.Range("F27").PasteSpecial(-4122) | Out-Null
$ws.Range("F30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F29").Clear() | Out-Null
$ws.Range("F32").Clear() | Out-Null

# ---------------------------------------------------------------------------
# Selection / scroll position reported by the saved view.
# ---------------------------------------------------------------------------

try {
    $excel.ActiveWindow.ScrollRow = 13
} catch {
}

$ws.Range("A2:A54").Select() | Out-Null

# ---------------------------------------------------------------------------
# Re-apply the (no-op, already-sorted) sort so the sheet's remembered
# sortState range stays in sync with the data.
# ---------------------------------------------------------------------------

$sortRange = $ws.Range("A2:D57")
$ws.Sort.SortFields.Clear() | Out-Null
$ws.Sort.SortFields.Add($ws.Range("C2:C57")) | Out-Null
$ws.Sort.SetRange($sortRange) | Out-Null
$ws.Sort.Header = -4142  # xlNo
$ws.Sort.Apply() | Out-Null
